$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "https://github.com/Akshat033692/IBM-PROJECT.git"
